$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new "BOUNDARY" header block (AS1:AZ1), mirrors the existing
#     8-column attack blocks (e.g. FGSM in AK1:AR1). Merge first while the
#     cells are still blank/default so the later format paste (from the
#     FGSM block) ends up applied uniformly to every cell in the merge. ---
$ws.Range("AS1:AZ1").Merge()
$ws.Range("AK1:AR1").Copy()
$ws.Range("AS1:AZ1").PasteSpecial(-4122)
$ws.Range("AS1").Value = "BOUNDARY"

# --- Row 2: epsilon labels stored as text ("0.01" ... "0.20"), matching
#     the existing row-2 style/type in the other attack blocks. ---
$ws.Range("AK2:AR2").Copy()
$ws.Range("AS2:AZ2").PasteSpecial(-4122)

$epsCols = @("AS","AT","AU","AV","AW","AX","AY","AZ")
$epsVals = @("0.01","0.02","0.03","0.04","0.05","0.07","0.10","0.20")
for ($i = 0; $i -lt $epsCols.Length; $i++) {
    $addr = $epsCols[$i] + "2"
    $ws.Range($addr).Formula = "=""" + $epsVals[$i] + """"
}
$ws.Range("AS2:AZ2").Copy()
$ws.Range("AS2:AZ2").PasteSpecial(-4163)

# --- Rows 4-12: numeric BOUNDARY-attack data (MAE / RMSE / SIM per model) ---
$ws.Range("AS4").Value = 4.338001728057861
$ws.Range("AT4").Value = 4.382240772247314
$ws.Range("AU4").Value = 4.409730911254883
$ws.Range("AV4").Value = 4.515683174133301
$ws.Range("AW4").Value = 4.662610054016113
$ws.Range("AX4").Value = 4.710359573364258
$ws.Range("AY4").Value = 4.960222244262695
$ws.Range("AZ4").Value = 6.929498672485352
$ws.Range("AS5").Value = 5.46575095640627
$ws.Range("AT5").Value = 5.490229684996662
$ws.Range("AU5").Value = 5.515694334118485
$ws.Range("AV5").Value = 5.671658430352267
$ws.Range("AW5").Value = 5.857137919300936
$ws.Range("AX5").Value = 5.960993815264784
$ws.Range("AY5").Value = 6.062538382811807
$ws.Range("AZ5").Value = 8.817718908427898
$ws.Range("AS6").Value = 0.9996315836906433
$ws.Range("AT6").Value = 0.9996266961097717
$ws.Range("AU6").Value = 0.9996223449707031
$ws.Range("AV6").Value = 0.9996002316474915
$ws.Range("AW6").Value = 0.9995725750923157
$ws.Range("AX6").Value = 0.9995550513267517
$ws.Range("AY6").Value = 0.9995250701904297
$ws.Range("AZ6").Value = 0.9989747405052185
$ws.Range("AS7").Value = 4.064916133880615
$ws.Range("AT7").Value = 4.169347763061523
$ws.Range("AU7").Value = 4.344825744628906
$ws.Range("AV7").Value = 4.57079029083252
$ws.Range("AW7").Value = 4.790963649749756
$ws.Range("AX7").Value = 5.10933256149292
$ws.Range("AY7").Value = 6.015357971191406
$ws.Range("AZ7").Value = 10.81198120117188
$ws.Range("AS8").Value = 5.060214810240891
$ws.Range("AT8").Value = 5.174336104262367
$ws.Range("AU8").Value = 5.335597014839178
$ws.Range("AV8").Value = 5.57603813051078
$ws.Range("AW8").Value = 5.933111235215749
$ws.Range("AX8").Value = 6.353197701371715
$ws.Range("AY8").Value = 7.634249970594714
$ws.Range("AZ8").Value = 13.28169519014886
$ws.Range("AS9").Value = 0.9997678995132446
$ws.Range("AT9").Value = 0.999756395816803
$ws.Range("AU9").Value = 0.9997288584709167
$ws.Range("AV9").Value = 0.9996874928474426
$ws.Range("AW9").Value = 0.9996163845062256
$ws.Range("AX9").Value = 0.9994912147521973
$ws.Range("AY9").Value = 0.9991987943649292
$ws.Range("AZ9").Value = 0.9972560405731201
$ws.Range("AS10").Value = 3.439592838287354
$ws.Range("AT10").Value = 3.569204568862915
$ws.Range("AU10").Value = 3.669927358627319
$ws.Range("AV10").Value = 4.01426887512207
$ws.Range("AW10").Value = 4.312092781066895
$ws.Range("AX10").Value = 5.090110301971436
$ws.Range("AY10").Value = 6.635757446289062
$ws.Range("AZ10").Value = 9.912277221679688
$ws.Range("AS11").Value = 4.306086761082471
$ws.Range("AT11").Value = 4.43080084801666
$ws.Range("AU11").Value = 4.606659480492764
$ws.Range("AV11").Value = 5.055466275014847
$ws.Range("AW11").Value = 5.355219261318064
$ws.Range("AX11").Value = 6.423956986013706
$ws.Range("AY11").Value = 8.181480466249914
$ws.Range("AZ11").Value = 12.60672214968566
$ws.Range("AS12").Value = 0.999796450138092
$ws.Range("AT12").Value = 0.9997882843017578
$ws.Range("AU12").Value = 0.9997532367706299
$ws.Range("AV12").Value = 0.9996988773345947
$ws.Range("AW12").Value = 0.9996252059936523
$ws.Range("AX12").Value = 0.9994310736656189
$ws.Range("AY12").Value = 0.9989506602287292
$ws.Range("AZ12").Value = 0.9974680542945862
